$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "carrier" column (D) for the practice rows (p2/p3/p4),
# mirroring the pair_kind already present in column K for those rows.
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Tag the generic-word practice pairs (rows 6-9) with which unique
# media type each pair needs.
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Populate pair_kind (C) / carrier (D) for rows 14-21, matching the
# unique_video / unique_audio grouping and carrier word pattern used
# above.
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"
$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"
$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"
$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"
$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"
$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"
$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"
$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
